$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.554.58'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.923.01'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "'326.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").Value = "'1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").Value = "'0.4069"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = "'0.08242"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'1.012"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").Value = "'23.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.86%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.900.71'
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'6.085"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").Value = "'7.282"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.28%  '
$ws.Range("D15").Value = "'91.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("D16").Value = "'0.06859"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").Value = "'1.012"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = "'17.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").Value = "'1.010"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("D21").Value = '29.566.22'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = "'5.685"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").Value = "'11.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = "'2.185"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").Value = '2.112.35'
$ws.Range("E25").Value = '  -2.07%  '
$ws.Range("D26").Value = "'155.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = "'6.489"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.54%  '
$ws.Range("D28").Value = "'20.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").Value = "'2.100"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = "'120.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("D31").Value = "'1.019"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("D32").Value = "'0.09661"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("D33").Value = "'5.635"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("D34").Value = "'3.554"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("D35").Value = "'1.376"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.15%  '
$ws.Range("D36").Value = "'0.06381"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.35%  '
$ws.Range("D37").Value = "'0.02300"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.30%  '
$ws.Range("D38").Value = "'1.188"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").Value = "'7.898"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("D42").Value = "'0.1853"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = "'2.481"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").Value = "'1.281"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").Value = "'12.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").Value = "'0.07509"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.97%  '
$ws.Range("D47").Value = "'0.5576"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").Value = "'1.953"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("D49").Value = "'119.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.94%  '
$ws.Range("D50").Value = "'2.437"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.51%  '
$ws.Range("D51").Value = "'72.23"
$ws.Range("D51").Style = "Normal"
